# Week 13 logging: add new players to WR and RB stat sheets.

$wb = $excel.ActiveWorkbook

# --- WR sheet: add K.Blanton as a new row (row 8) ---
$wr = $wb.Worksheets.Item("WR")
$wr.Cells.Item(8, 1).Value = "K.Blanton"
for ($c = 2; $c -le 10; $c++) {
    $wr.Cells.Item(8, $c).Value = 0
}
$wr.Range("I9").Select()

# --- RB sheet: add M.Sargent as a new row (row 6) ---
$rb = $wb.Worksheets.Item("RB")
$rb.Cells.Item(6, 1).Value = "M.Sargent"
for ($c = 2; $c -le 10; $c++) {
    $rb.Cells.Item(6, $c).Value = 0
}
$rb.Range("J7").Select()

# RB becomes the active/selected tab after this logging session.
$rb.Activate()
$rb.Range("J7").Select()
